# Split the single "test" run into three runs: "T" / "est" / " change",
# matching the target OOXML diff (same text, no real formatting change,
# just new run boundaries at index 1 and 4).

$d = $word.ActiveDocument

# --- "test" -> "T" + "est" -------------------------------------------------
# Capitalize just the leading "t" in place.
$rT = $d.Range(0, 1)
$rT.Text = "T"

# A plain text edit on adjacent, identically-formatted runs gets coalesced
# back into a single run when the document is saved. Toggling a character
# format on/off (leaving the net formatting unchanged) is enough to pin the
# run boundary between "T" and "est" so the two stay as separate <w:r>
# elements.
$rT.Bold = $true
$rT.Bold = $false
Write-Output "After capitalizing: $($d.Content.Text)"

# --- append " change" as its own run ---------------------------------------
$rTail = $d.Range(4, 4)
$rTail.InsertAfter(" change")
Write-Output "Final text: $($d.Content.Text)"
